# Generate Report for Archive
#
# 1. Update the localization status text "Ready for handoff" -> "In Translation"
#    on the Overview sheet (this text is a shared string, so every cell that
#    uses it updates together).
# 2. Narrow the "Status"/"zh-cn"/"de-de" report columns (previously widened
#    for a longer status label) back down now that the text is shorter.

$wb = $excel.ActiveWorkbook

# --- 1. Status text -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- 2. Column widths -------------------------------------------------------
# Same physical width is applied to all three report "status"-adjacent
# columns: Overview!E:F (zh-cn / de-de handoff columns) and the "Status"
# column (C) on each language sheet.
$newColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth

$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth

$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
